$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(59).Insert()
$ws.Cells.Item(59,1).Value = 10
$ws.Cells.Item(59,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(59,3).Value = "La Araucanía"
$ws.Cells.Item(59,4).Value = 45174
$ws.Cells.Item(59,5).Value = 9
$ws.Cells.Item(59,6).Value = "Fruta"
$ws.Cells.Item(59,7).Value = 100108
$ws.Cells.Item(59,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(59,9).Value = 100108004
$ws.Cells.Item(59,10).Value = "Papaya"
$ws.Cells.Item(59,11).Value = "Cultivar IV Región"
$ws.Cells.Item(59,12).Value = "Primera"
$ws.Cells.Item(59,13).Value = 80
$ws.Cells.Item(59,14).Value = 24000
$ws.Cells.Item(59,15).Value = 24000
$ws.Cells.Item(59,16).Value = 24000
$ws.Cells.Item(59,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(59,18).Value = "Provincia del Elquí"
$ws.Cells.Item(59,19).Value = 2400
$ws.Cells.Item(59,20).Value = 10
